$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value = 11
$ws.Cells.Item(3, 5).Value = 11
$ws.Cells.Item(4, 5).Value = 11
$ws.Cells.Item(5, 5).Value = 5
$ws.Cells.Item(6, 5).Value = 11
$ws.Cells.Item(7, 5).Value = 5
$ws.Cells.Item(8, 5).Value = 11
$ws.Cells.Item(9, 5).Value = 5
$ws.Cells.Item(10, 5).Value = 4
$ws.Cells.Item(11, 5).Value = 11
$ws.Cells.Item(12, 5).Value = 5
$ws.Cells.Item(13, 5).Value = 11
$ws.Cells.Item(14, 5).Value = 11
$ws.Cells.Item(15, 5).Value = 11
$ws.Cells.Item(16, 5).Value = 9
$ws.Cells.Item(17, 5).Value = 5
$ws.Cells.Item(18, 5).Value = 8
$ws.Cells.Item(19, 5).Value = 8
$ws.Cells.Item(20, 5).Value = 8
$ws.Cells.Item(21, 5).Value = 8
$ws.Cells.Item(22, 5).Value = 5
$ws.Cells.Item(23, 5).Value = 5
$ws.Cells.Item(24, 5).Value = 5
$ws.Cells.Item(25, 5).Value = 5
$ws.Cells.Item(26, 5).Value = 5
$ws.Cells.Item(27, 5).Value = 5
$ws.Cells.Item(28, 5).Value = 8
$ws.Cells.Item(29, 5).Value = 8
$ws.Cells.Item(30, 5).Value = 8
$ws.Cells.Item(31, 5).Value = 8
$ws.Cells.Item(32, 5).Value = 8
$ws.Cells.Item(33, 5).Value = 8
$ws.Cells.Item(34, 5).Value = 8
$ws.Cells.Item(35, 5).Value = 8
$ws.Cells.Item(37, 5).Value = 8
$ws.Cells.Item(38, 5).Value = 8
$ws.Cells.Item(39, 5).Value = 8
$ws.Cells.Item(40, 5).Value = 4
$ws.Cells.Item(41, 5).Value = 4
$ws.Cells.Item(42, 5).Value = 8
$ws.Cells.Item(43, 5).Value = 5
$ws.Cells.Item(44, 5).Value = 4
$ws.Cells.Item(45, 5).Value = 5
$ws.Cells.Item(46, 5).Value = 4
$ws.Cells.Item(47, 5).Value = 8
$ws.Cells.Item(48, 5).Value = 4
$ws.Cells.Item(49, 5).Value = 5
$ws.Cells.Item(50, 5).Value = 3
$ws.Cells.Item(51, 5).Value = 3
$ws.Cells.Item(52, 5).Value = 3
$ws.Cells.Item(53, 5).Value = 3
$ws.Cells.Item(54, 5).Value = 3
$ws.Cells.Item(55, 5).Value = 3
$ws.Cells.Item(56, 5).Value = 3
$ws.Cells.Item(57, 5).Value = 3
$ws.Cells.Item(58, 5).Value = 7
$ws.Cells.Item(59, 5).Value = 7
$ws.Cells.Item(60, 5).Value = 7
$ws.Cells.Item(61, 5).Value = 5
$ws.Cells.Item(62, 5).Value = 7
$ws.Cells.Item(63, 5).Value = 7
$ws.Cells.Item(64, 5).Value = 7
$ws.Cells.Item(65, 5).Value = 8
$ws.Cells.Item(66, 5).Value = 8
$ws.Cells.Item(67, 5).Value = 8
$ws.Cells.Item(68, 5).Value = 8
$ws.Cells.Item(69, 5).Value = 8
$ws.Cells.Item(70, 5).Value = 9
$ws.Cells.Item(71, 5).Value = 9
$ws.Cells.Item(72, 5).Value = 9
$ws.Cells.Item(73, 5).Value = 9
$ws.Cells.Item(74, 5).Value = 9
$ws.Cells.Item(75, 5).Value = 9
$ws.Cells.Item(76, 5).Value = 9
$ws.Cells.Item(77, 5).Value = 2
$ws.Cells.Item(78, 5).Value = 2
$ws.Cells.Item(79, 5).Value = 2
$ws.Cells.Item(80, 5).Value = 2
$ws.Cells.Item(81, 5).Value = 2
$ws.Cells.Item(82, 5).Value = 2
$ws.Cells.Item(83, 5).Value = 2
$ws.Cells.Item(84, 5).Value = 2
$ws.Cells.Item(85, 5).Value = 2
$ws.Cells.Item(86, 5).Value = 2
$ws.Cells.Item(87, 5).Value = 4
$ws.Cells.Item(88, 5).Value = 4
$ws.Cells.Item(89, 5).Value = 4
$ws.Cells.Item(90, 5).Value = 4
$ws.Cells.Item(91, 5).Value = 5
$ws.Cells.Item(92, 5).Value = 4
$ws.Cells.Item(93, 5).Value = 2
$ws.Cells.Item(94, 5).Value = 7
$ws.Cells.Item(94, 6).Value = 20260129
$ws.Cells.Item(95, 5).Value = 1
$ws.Cells.Item(96, 5).Value = 9
$ws.Cells.Item(97, 5).Value = 9
$ws.Cells.Item(98, 5).Value = 9
$ws.Cells.Item(99, 5).Value = 9
